# Delete row 9 (Equipment Number GLDU9669533 / WO 9071905279-01 / BOL
# ONEYSZPV37250300) from the data table on Sheet1.  Excel shifts every
# row below it up by one and the shared-strings table is recompacted,
# dropping the three string values that were only referenced by the
# deleted row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(9).Delete()
